# Edit: populate the "Testing Table for Ordering Pizza Combos" worksheet
# (Selenium-results grid for FF/IE/Chrome + a device checklist) on Sheet1,
# rows 5-13. Cell writes are ordered so that shared-string indices land in
# the same sequence as the target workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range('B5').Value = 'Yes'
$ws.Range('A5').Value = 'FF'
$ws.Range('A6').Value = 'IE'
$ws.Range('A7').Value = 'Chrome'
$ws.Range('D6').Value = 'No--can''t remove pizza. Image off plum. Price wrong: £7.94 instead of £7.49'
$ws.Range('B6').Value = 'No--can''t remove pizza. Image off plum. Says "Classic Deluxe" instead of "Classic Deluxe Pizza". Price fine.'
$ws.Range('C6').Value = 'No--can''t remove pizza. Image off plum. Price fine.'
$ws.Range('E6').Value = 'No--can''t remove pizza. Price fine.'
$ws.Range('F6').Value = 'No--can''t remove either pizza. Intended price is wrong: should be 16.48, not 16.93'
$ws.Range('G6').Value = 'No--can remove neither pizza. Price wrong: should be £13.48'
$ws.Range('H6').Value = 'No--can''t remove any pizza. Price should be £22.47, not 22.92'
$ws.Range('I6').Value = 'No--can''t remove any pizza. Price fine.'
$ws.Range('J6').Value = 'No--can''t remove any pizza. Price wrong: should be 25.47, not 25.92'
$ws.Range('K6').Value = 'Price fine; can''t remove any pizza.'
$ws.Range('L6').Value = 'No--price should be 23.97 not 24.87.'
$ws.Range('M6').Value = 'No--can''t remove pizzas; price should be 20.97, not 21.87'
$ws.Range('O6').Value = 'Price fine; can''t remove pizzas, though.'
$ws.Range('P6').Value = 'No--can''t remove pizzas; price should be 32.96, not 33.86'
$ws.Range('Q6').Value = 'No--price should be 26.96, not 27.86; can''t remove any pizza.'
$ws.Range('R6').Value = 'No--price should be 37.45, not 37.9; can''t remove pizzas.'
$ws.Range('D7').Value = 'No--pizza should be 7.49, not 7.94'
$ws.Range('F7').Value = 'No--should be 16.48, not 16.93'
$ws.Range('G5').Value = 'No--should be 13.48, not 13.93'
$ws.Range('H5').Value = 'No--should be 22.47, not 22.92'
$ws.Range('J5').Value = 'No--should be 25.47, not 25.92'
$ws.Range('L5').Value = 'No--price should be 23.97, not 24.87'
$ws.Range('M5').Value = 'No--price should be 20.97, not 21.87'
$ws.Range('N5').Value = 'No--price should be 19.47, not 19.92'
$ws.Range('P5').Value = 'No--price should be 32.96, not 33.86'
$ws.Range('Q5').Value = 'No--should be 26.96, not 27.86'
$ws.Range('R5').Value = 'No--price should be 37.45, not 37.9'
$ws.Range('S5').Value = 'No--price should be 38.95, not 39.85'
$ws.Range('T5').Value = 'No--price should be 35.95, not 36.85'
$ws.Range('U5').Value = 'No--price should be 44.94, not 45.84'
$ws.Range('D5').Value = 'No--price should be 7.49, not 7.94'
$ws.Range('S6').Value = 'No--price should be 38.95, not 39.85. Can''t remove any pizza.'
$ws.Range('T6').Value = 'No--price should be 35.95, not 36.85. Can''t remove any pizza.'
$ws.Range('U6').Value = 'No--price should be 44.94, not 45.84. Can''t remove pizzas.'
$ws.Range('F5').Value = 'No--price should be 16.48, not 16.93.'
$ws.Range('N6').Value = 'No--can''t remove pizzas; price should be 19.47, not 19.92'
$ws.Range('A9').Value = 'Kindle Fire'
$ws.Range('A10').Value = 'iPad4'
$ws.Range('A11').Value = 'iPhone 6'
$ws.Range('A12').Value = 'Nexus 10'
$ws.Range('A13').Value = 'Galaxy S4'
$ws.Range('B11').Value = 'Images off plum. Price fine.'
$ws.Range('B13').Value = 'Images off plum. Price fine. '
$ws.Range('D13').Value = 'No--pizza should be 7.49, not 7.94. Images off plum. '
$ws.Range('F13').Value = 'No--should be 16.48, not 16.93. Images off plum.'
$ws.Range('G13').Value = 'No--should be 13.48, not 13.93. Images off plum.'
$ws.Range('J13').Value = 'No--should be 25.47, not 25.92. Images off plum.'

# Remaining cells repeat strings already introduced above.
$ws.Range('C5').Value = 'Yes'
$ws.Range('E5').Value = 'Yes'
$ws.Range('I5').Value = 'Yes'
$ws.Range('K5').Value = 'Yes'
$ws.Range('O5').Value = 'Yes'
$ws.Range('B7').Value = 'Yes'
$ws.Range('C7').Value = 'Yes'
$ws.Range('E7').Value = 'Yes'
$ws.Range('G7').Value = 'No--should be 13.48, not 13.93'
$ws.Range('H7').Value = 'No--should be 22.47, not 22.92'
$ws.Range('I7').Value = 'Yes'
$ws.Range('J7').Value = 'No--should be 25.47, not 25.92'
$ws.Range('K7').Value = 'Yes'
$ws.Range('L7').Value = 'No--price should be 23.97, not 24.87'
$ws.Range('M7').Value = 'No--price should be 20.97, not 21.87'
$ws.Range('N7').Value = 'No--price should be 19.47, not 19.92'
$ws.Range('O7').Value = 'Yes'
$ws.Range('P7').Value = 'No--price should be 32.96, not 33.86'
$ws.Range('Q7').Value = 'No--should be 26.96, not 27.86'
$ws.Range('R7').Value = 'No--price should be 37.45, not 37.9'
$ws.Range('S7').Value = 'No--price should be 38.95, not 39.85'
$ws.Range('T7').Value = 'No--price should be 35.95, not 36.85'
$ws.Range('U7').Value = 'No--price should be 44.94, not 45.84'
$ws.Range('B9').Value = 'Yes'
$ws.Range('B10').Value = 'Yes'
$ws.Range('B12').Value = 'Yes'
$ws.Range('C13').Value = 'Images off plum. Price fine. '
$ws.Range('H13').Value = 'No--should be 22.47, not 22.92'
$ws.Range('I13').Value = 'Images off plum. Price fine.'

# Leave the selection where the author left off.
$ws.Range('K13').Select() | Out-Null
